# Reorders the "Recorded By" (column G) value in each data row so that the
# entry "System" (exact case) is moved to the front of the comma-separated
# list, preserving the relative order of the remaining entries.
#
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#      "system, System, backup@backdoor.com" -> "System, system, backup@backdoor.com"
#
# Rows where "System" is not present, or is already first, or where it is
# the only entry, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$changed = 0

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notlike "*,*") { continue }

    $rawParts = $val -split ","
    $parts = @()
    foreach ($p in $rawParts) { $parts += $p.Trim() }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    if ($hasSystem -and -not $parts[0].Equals("System")) {
        $newParts = @("System")
        foreach ($p in $parts) {
            if (-not $p.Equals("System")) { $newParts += $p }
        }
        $newVal = [string]::Join(", ", $newParts)

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
            $changed++
        }
    }
}

Write-Host "Updated $changed cells in column G"
